$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update corrected values on existing rows ---

# Row 267 (2025-10-08)
$ws.Cells.Item(267, 2).Value = 247875.17
$ws.Cells.Item(267, 3).Value = 3990.05
$ws.Cells.Item(267, 4).Value = 485.76
$ws.Cells.Item(267, 5).Value = 251379.46
$ws.Cells.Item(267, 10).Value = 352987.37

# Row 268 (2025-10-09)
$ws.Cells.Item(268, 2).Value = 278024.34
$ws.Cells.Item(268, 4).Value = 4052.13
$ws.Cells.Item(268, 5).Value = 277269.52
$ws.Cells.Item(268, 10).Value = 333295.39

# Row 274 (2025-10-15)
$ws.Cells.Item(274, 2).Value = 95471.50999999999
$ws.Cells.Item(274, 4).Value = 612.92
$ws.Cells.Item(274, 5).Value = 96805.09
$ws.Cells.Item(274, 10).Value = 119826.94

# Row 276 (2025-10-17)
$ws.Cells.Item(276, 2).Value = 89316.2
$ws.Cells.Item(276, 4).Value = 76.73999999999999
$ws.Cells.Item(276, 5).Value = 89842.49999999999
$ws.Cells.Item(276, 10).Value = 106043.07

# --- Append new rows 279-285 ---

# Row 279 (2025-10-20)
$ws.Cells.Item(279, 1).Value = 45950
$ws.Cells.Item(279, 2).Value = 63529.98
$ws.Cells.Item(279, 3).Value = 2685.74
$ws.Cells.Item(279, 4).Value = 16029.71
$ws.Cells.Item(279, 5).Value = 50186.01
$ws.Cells.Item(279, 6).Value = 25973.77
$ws.Cells.Item(279, 7).Value = 0
$ws.Cells.Item(279, 8).Value = 2000
$ws.Cells.Item(279, 9).Value = 23973.77
$ws.Cells.Item(279, 10).Value = 74159.78
$ws.Cells.Item(279, 11).Value = "SOJA"

# Row 280 (2025-10-21)
$ws.Cells.Item(280, 1).Value = 45951
$ws.Cells.Item(280, 2).Value = 68171.38
$ws.Cells.Item(280, 3).Value = 364.01
$ws.Cells.Item(280, 4).Value = 6490
$ws.Cells.Item(280, 5).Value = 62045.39
$ws.Cells.Item(280, 6).Value = 9833.360000000001
$ws.Cells.Item(280, 7).Value = 0
$ws.Cells.Item(280, 8).Value = 0
$ws.Cells.Item(280, 9).Value = 9833.360000000001
$ws.Cells.Item(280, 10).Value = 71878.75
$ws.Cells.Item(280, 11).Value = "SOJA"

# Row 281 (2025-10-22)
$ws.Cells.Item(281, 1).Value = 45952
$ws.Cells.Item(281, 2).Value = 51425.44
$ws.Cells.Item(281, 3).Value = 728.95
$ws.Cells.Item(281, 4).Value = 10.86
$ws.Cells.Item(281, 5).Value = 52143.53
$ws.Cells.Item(281, 6).Value = 10941.52
$ws.Cells.Item(281, 7).Value = 0
$ws.Cells.Item(281, 8).Value = 0
$ws.Cells.Item(281, 9).Value = 10941.52
$ws.Cells.Item(281, 10).Value = 63085.05
$ws.Cells.Item(281, 11).Value = "SOJA"

# Row 282 (2025-10-23)
$ws.Cells.Item(282, 1).Value = 45953
$ws.Cells.Item(282, 2).Value = 60365.42
$ws.Cells.Item(282, 3).Value = 644.3199999999999
$ws.Cells.Item(282, 4).Value = 2000
$ws.Cells.Item(282, 5).Value = 59009.74
$ws.Cells.Item(282, 6).Value = 14164.01
$ws.Cells.Item(282, 7).Value = 248.65
$ws.Cells.Item(282, 8).Value = 0
$ws.Cells.Item(282, 9).Value = 14412.66
$ws.Cells.Item(282, 10).Value = 73422.39999999999
$ws.Cells.Item(282, 11).Value = "SOJA"

# Row 283 (2025-10-24)
$ws.Cells.Item(283, 1).Value = 45954
$ws.Cells.Item(283, 2).Value = 72460.7
$ws.Cells.Item(283, 3).Value = 297.46
$ws.Cells.Item(283, 4).Value = 1060
$ws.Cells.Item(283, 5).Value = 71698.16
$ws.Cells.Item(283, 6).Value = 54767.44
$ws.Cells.Item(283, 7).Value = 90
$ws.Cells.Item(283, 8).Value = 2928
$ws.Cells.Item(283, 9).Value = 51929.44
$ws.Cells.Item(283, 10).Value = 123627.6
$ws.Cells.Item(283, 11).Value = "SOJA"

# Row 284 (2025-10-25)
$ws.Cells.Item(284, 1).Value = 45955
$ws.Cells.Item(284, 2).Value = 61.36
$ws.Cells.Item(284, 3).Value = 0
$ws.Cells.Item(284, 4).Value = 0
$ws.Cells.Item(284, 5).Value = 61.36
$ws.Cells.Item(284, 6).Value = 0
$ws.Cells.Item(284, 7).Value = 0
$ws.Cells.Item(284, 8).Value = 0
$ws.Cells.Item(284, 9).Value = 0
$ws.Cells.Item(284, 10).Value = 61.36
$ws.Cells.Item(284, 11).Value = "SOJA"

# Row 285 (2025-10-27)
$ws.Cells.Item(285, 1).Value = 45957
$ws.Cells.Item(285, 2).Value = 2677.84
$ws.Cells.Item(285, 3).Value = 0
$ws.Cells.Item(285, 4).Value = 0
$ws.Cells.Item(285, 5).Value = 2677.84
$ws.Cells.Item(285, 6).Value = 0
$ws.Cells.Item(285, 7).Value = 0
$ws.Cells.Item(285, 8).Value = 0
$ws.Cells.Item(285, 9).Value = 0
$ws.Cells.Item(285, 10).Value = 2677.84
$ws.Cells.Item(285, 11).Value = "SOJA"

# Ensure column A date style (yyyy-mm-dd, style index 2) is applied for the new rows,
# matching the existing column A formatting.
$ws.Range("A279:A285").NumberFormat = "yyyy-mm-dd"
